$d = $word.ActiveDocument

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14 = "http://schemas.microsoft.com/office/word/2010/wordml"

# Locate the paragraph "Commit 4:" that sits right under "JS2.js contains:"
# (there are two "Commit 4:" paragraphs in the document; this is the one
# whose next sibling is the trailing blank paragraph at the very end of
# the body).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Trim() -eq "Commit 4:") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs($targetIndex)

# Rebuild the "Commit 4:" paragraph in place: drop the paragraph-mark's
# <w:lang> and append a new run holding a single trailing space (also
# without an explicit <w:lang>), while keeping the paragraph's own
# identity attributes untouched.
$target.Range.InsertXML('<w:p xmlns:w="' + $w + '" xmlns:w14="' + $w14 + '" w14:paraId="39882B65" w14:textId="597257D5" w:rsidR="00A20BD5" w:rsidRPr="00A20BD5" w:rsidRDefault="00A20BD5" w:rsidP="00A20BD5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Commit 4:</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>')

# Append a new (currently empty) paragraph right after it, then fill it
# in with the "While loop ..." text.
$target = $d.Paragraphs($targetIndex)
$target.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($targetIndex + 1)
$p1.Range.InsertXML('<w:p xmlns:w="' + $w + '"><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>While loop, for loop, scope of  variable using let and var declarations, break, continue, do while loop</w:t></w:r></w:p>')

# Append a new paragraph for "Commit 5:".
$p1 = $d.Paragraphs($targetIndex + 1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($targetIndex + 2)
$p2.Range.InsertXML('<w:p xmlns:w="' + $w + '"><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Commit 5:</w:t></w:r></w:p>')

# Append a new paragraph for the "Arrays, ..." text.
$p2 = $d.Paragraphs($targetIndex + 2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($targetIndex + 3)
$p3.Range.InsertXML('<w:p xmlns:w="' + $w + '"><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Arrays, push , pop, shift, unshift, primitive vs reference types, array cloning methods, creating arrays out of already existing arrays with additional items, array concatenation</w:t></w:r></w:p>')
